# Auto-applies the updated TPM-derived values for Efnb1-Ephb3 LR pairs (rows 2-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 10.430598
$ws.Range("H2").Value2 = 20.861196
$ws.Range("I2").Value2 = 0.5361372314060847
$ws.Range("J2").Value2 = 0.4628062232731417
$ws.Range("M2").Value2 = 0.3250655
$ws.Range("N2").Value2 = 0.650131
$ws.Range("O2").Value2 = 0.03171104771407953
$ws.Range("P2").Value2 = 0.02364479350543662
$ws.Range("Q2").Value2 = 3.390627554169
$ws.Range("R2").Value2 = 13.562510216676
$ws.Range("S2").Value2 = 0.01700147332641285
$ws.Range("T2").Value2 = 0.01094295758232443

# Row 3
$ws.Range("G3").Value2 = 10.430598
$ws.Range("H3").Value2 = 20.861196
$ws.Range("I3").Value2 = 0.5361372314060847
$ws.Range("J3").Value2 = 0.4628062232731417
$ws.Range("O3").Value2 = 0.6816872822276142
$ws.Range("P3").Value2 = 0.762432473166021
$ws.Range("Q3").Value2 = 72.887774105972
$ws.Range("R3").Value2 = 437.326644635832
$ws.Range("S3").Value2 = 0.3654779321782514
$ws.Range("T3").Value2 = 0.3528584934067672

# Row 4
$ws.Range("G4").Value2 = 10.430598
$ws.Range("H4").Value2 = 20.861196
$ws.Range("I4").Value2 = 0.5361372314060847
$ws.Range("J4").Value2 = 0.4628062232731417
$ws.Range("M4").Value2 = 2.9317775
$ws.Range("N4").Value2 = 5.863555
$ws.Range("O4").Value2 = 0.2860030861151516
$ws.Range("P4").Value2 = 0.2132532477035712
$ws.Range("Q4").Value2 = 30.580192527945
$ws.Range("R4").Value2 = 122.32077011178
$ws.Range("S4").Value2 = 0.1533369027633734
$ws.Range("T4").Value2 = 0.09869493017042157

# Row 5
$ws.Range("D5").Value2 = "Resolving-Mac"
$ws.Range("G5").Value2 = 10.430598
$ws.Range("H5").Value2 = 20.861196
$ws.Range("I5").Value2 = 0.5361372314060847
$ws.Range("J5").Value2 = 0.4628062232731417
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.006136
$ws.Range("N5").Value2 = 0.018408
$ws.Range("O5").Value2 = 0.0005985839431548165
$ws.Range("P5").Value2 = 0.0006694856249710864
$ws.Range("Q5").Value2 = 0.064002149328
$ws.Range("R5").Value2 = 0.384012895968
$ws.Range("S5").Value2 = 0.0003209231380471605
$ws.Range("T5").Value2 = 0.0003098421136285275

# Row 6
$ws.Range("I6").Value2 = 0.2610112166105907
$ws.Range("J6").Value2 = 0.3379664990099709
$ws.Range("M6").Value2 = 0.3250655
$ws.Range("N6").Value2 = 0.650131
$ws.Range("O6").Value2 = 0.03171104771407953
$ws.Range("P6").Value2 = 0.02364479350543662
$ws.Range("Q6").Value2 = 1.650681525448333
$ws.Range("R6").Value2 = 9.90408915269
$ws.Range("S6").Value2 = 0.008276939143848389
$ws.Range("T6").Value2 = 0.007991148080846113

# Row 7
$ws.Range("I7").Value2 = 0.2610112166105907
$ws.Range("J7").Value2 = 0.3379664990099709
$ws.Range("O7").Value2 = 0.6816872822276142
$ws.Range("P7").Value2 = 0.762432473166021
$ws.Range("S7").Value2 = 0.1779280268821967
$ws.Range("T7").Value2 = 0.2576766336874337

# Row 8
$ws.Range("I8").Value2 = 0.2610112166105907
$ws.Range("J8").Value2 = 0.3379664990099709
$ws.Range("M8").Value2 = 2.9317775
$ws.Range("N8").Value2 = 5.863555
$ws.Range("O8").Value2 = 0.2860030861151516
$ws.Range("P8").Value2 = 0.2132532477035712
$ws.Range("Q8").Value2 = 14.88755637240833
$ws.Range("R8").Value2 = 89.32533823444999
$ws.Range("S8").Value2 = 0.07465001346129925
$ws.Range("T8").Value2 = 0.07207245352888207

# Row 9
$ws.Range("D9").Value2 = "Resolving-Mac"
$ws.Range("I9").Value2 = 0.2610112166105907
$ws.Range("J9").Value2 = 0.3379664990099709
$ws.Range("K9").Value2 = 1
$ws.Range("L9").Value2 = 0.3333333333333333
$ws.Range("M9").Value2 = 0.006136
$ws.Range("N9").Value2 = 0.018408
$ws.Range("O9").Value2 = 0.0005985839431548165
$ws.Range("P9").Value2 = 0.0006694856249710864
$ws.Range("Q9").Value2 = 0.03115858754666666
$ws.Range("R9").Value2 = 0.28042728792
$ws.Range("S9").Value2 = 0.0001562371232464033
$ws.Range("T9").Value2 = 0.0002262637128089804

# Row 10
$ws.Range("G10").Value2 = 0.143788
$ws.Range("H10").Value2 = 0.431364
$ws.Range("I10").Value2 = 0.007390765153581618
$ws.Range("J10").Value2 = 0.009569822540183962
$ws.Range("M10").Value2 = 0.3250655
$ws.Range("N10").Value2 = 0.650131
$ws.Range("O10").Value2 = 0.03171104771407953
$ws.Range("P10").Value2 = 0.02364479350543662
$ws.Range("Q10").Value2 = 0.046740518114
$ws.Range("R10").Value2 = 0.280443108684
$ws.Range("S10").Value2 = 0.000234368906428783
$ws.Range("T10").Value2 = 0.0002262764778463228

# Row 11
$ws.Range("G11").Value2 = 0.143788
$ws.Range("H11").Value2 = 0.431364
$ws.Range("I11").Value2 = 0.007390765153581618
$ws.Range("J11").Value2 = 0.009569822540183962
$ws.Range("O11").Value2 = 0.6816872822276142
$ws.Range("P11").Value2 = 0.762432473166021
$ws.Range("Q11").Value2 = 1.004773385298667
$ws.Range("R11").Value2 = 9.042960467687999
$ws.Range("S11").Value2 = 0.005038190611127608
$ws.Range("T11").Value2 = 0.007296343467072392

# Row 12
$ws.Range("G12").Value2 = 0.143788
$ws.Range("H12").Value2 = 0.431364
$ws.Range("I12").Value2 = 0.007390765153581618
$ws.Range("J12").Value2 = 0.009569822540183962
$ws.Range("M12").Value2 = 2.9317775
$ws.Range("N12").Value2 = 5.863555
$ws.Range("O12").Value2 = 0.2860030861151516
$ws.Range("P12").Value2 = 0.2132532477035712
$ws.Range("Q12").Value2 = 0.42155442317
$ws.Range("R12").Value2 = 2.52932653902
$ws.Range("S12").Value2 = 0.002113781642676665
$ws.Range("T12").Value2 = 0.002040795736641069

# Row 13
$ws.Range("D13").Value2 = "Resolving-Mac"
$ws.Range("G13").Value2 = 0.143788
$ws.Range("H13").Value2 = 0.431364
$ws.Range("I13").Value2 = 0.007390765153581618
$ws.Range("J13").Value2 = 0.009569822540183962
$ws.Range("K13").Value2 = 1
$ws.Range("L13").Value2 = 0.3333333333333333
$ws.Range("M13").Value2 = 0.006136
$ws.Range("N13").Value2 = 0.018408
$ws.Range("O13").Value2 = 0.0005985839431548165
$ws.Range("P13").Value2 = 0.0006694856249710864
$ws.Range("Q13").Value2 = 0.000882283168
$ws.Range("R13").Value2 = 0.007940548512
$ws.Range("S13").Value2 = 0.000004423993348562097
$ws.Range("T13").Value2 = 0.000006406858624177449

# Row 14
$ws.Range("G14").Value2 = 2.8592275
$ws.Range("H14").Value2 = 5.718455000000001
$ws.Range("I14").Value2 = 0.1469655254483148
$ws.Range("J14").Value2 = 0.1268640859089486
$ws.Range("M14").Value2 = 0.3250655
$ws.Range("N14").Value2 = 0.650131
$ws.Range("O14").Value2 = 0.03171104771407953
$ws.Range("P14").Value2 = 0.02364479350543662
$ws.Range("Q14").Value2 = 0.9294362169012501
$ws.Range("R14").Value2 = 3.717744867605
$ws.Range("S14").Value2 = 0.004660430789816279
$ws.Range("T14").Value2 = 0.002999675114573061

# Row 15
$ws.Range("G15").Value2 = 2.8592275
$ws.Range("H15").Value2 = 5.718455000000001
$ws.Range("I15").Value2 = 0.1469655254483148
$ws.Range("J15").Value2 = 0.1268640859089486
$ws.Range("O15").Value2 = 0.6816872822276142
$ws.Range("P15").Value2 = 0.762432473166021
$ws.Range("Q15").Value2 = 19.97994056885167
$ws.Range("R15").Value2 = 119.87964341311
$ws.Range("S15").Value2 = 0.100184529624015
$ws.Range("T15").Value2 = 0.09672529877550622

# Row 16
$ws.Range("G16").Value2 = 2.8592275
$ws.Range("H16").Value2 = 5.718455000000001
$ws.Range("I16").Value2 = 0.1469655254483148
$ws.Range("J16").Value2 = 0.1268640859089486
$ws.Range("M16").Value2 = 2.9317775
$ws.Range("N16").Value2 = 5.863555
$ws.Range("O16").Value2 = 0.2860030861151516
$ws.Range("P16").Value2 = 0.2132532477035712
$ws.Range("Q16").Value2 = 8.382618851881251
$ws.Range("R16").Value2 = 33.530475407525
$ws.Range("S16").Value2 = 0.04203259383075288
$ws.Range("T16").Value2 = 0.02705417833702814

# Row 17
$ws.Range("D17").Value2 = "Resolving-Mac"
$ws.Range("G17").Value2 = 2.8592275
$ws.Range("H17").Value2 = 5.718455000000001
$ws.Range("I17").Value2 = 0.1469655254483148
$ws.Range("J17").Value2 = 0.1268640859089486
$ws.Range("K17").Value2 = 1
$ws.Range("L17").Value2 = 0.3333333333333333
$ws.Range("M17").Value2 = 0.006136
$ws.Range("N17").Value2 = 0.018408
$ws.Range("O17").Value2 = 0.0005985839431548165
$ws.Range("P17").Value2 = 0.0006694856249710864
$ws.Range("Q17").Value2 = 0.01754421994
$ws.Range("R17").Value2 = 0.10526531964
$ws.Range("S17").Value2 = 0.00008797120373067179
$ws.Range("T17").Value2 = 0.00008493368184113803

# Row 18
$ws.Range("G18").Value2 = 0.7979563333333332
$ws.Range("H18").Value2 = 2.393869
$ws.Range("I18").Value2 = 0.04101529934681446
$ws.Range("J18").Value2 = 0.05310805147033049
$ws.Range("M18").Value2 = 0.3250655
$ws.Range("N18").Value2 = 0.650131
$ws.Range("O18").Value2 = 0.03171104771407953
$ws.Range("P18").Value2 = 0.02364479350543662
$ws.Range("Q18").Value2 = 0.2593880744731666
$ws.Range("R18").Value2 = 1.556328446839
$ws.Range("S18").Value2 = 0.001300638114594088
$ws.Range("T18").Value2 = 0.001255728910492064

# Row 19
$ws.Range("G19").Value2 = 0.7979563333333332
$ws.Range("H19").Value2 = 2.393869
$ws.Range("I19").Value2 = 0.04101529934681446
$ws.Range("J19").Value2 = 0.05310805147033049
$ws.Range("O19").Value2 = 0.6816872822276142
$ws.Range("P19").Value2 = 0.762432473166021
$ws.Range("Q19").Value2 = 5.576023634544221
$ws.Range("R19").Value2 = 50.18421271089799
$ws.Range("S19").Value2 = 0.02795960794148199
$ws.Range("T19").Value2 = 0.04049130302755241

# Row 20
$ws.Range("G20").Value2 = 0.7979563333333332
$ws.Range("H20").Value2 = 2.393869
$ws.Range("I20").Value2 = 0.04101529934681446
$ws.Range("J20").Value2 = 0.05310805147033049
$ws.Range("M20").Value2 = 2.9317775
$ws.Range("N20").Value2 = 5.863555
$ws.Range("O20").Value2 = 0.2860030861151516
$ws.Range("P20").Value2 = 0.2132532477035712
$ws.Range("Q20").Value2 = 2.339430424049166
$ws.Range("R20").Value2 = 14.036582544295
$ws.Range("S20").Value2 = 0.0117305021911257
$ws.Range("T20").Value2 = 0.0113254644552564

# Row 21
$ws.Range("D21").Value2 = "Resolving-Mac"
$ws.Range("G21").Value2 = 0.7979563333333332
$ws.Range("H21").Value2 = 2.393869
$ws.Range("I21").Value2 = 0.04101529934681446
$ws.Range("J21").Value2 = 0.05310805147033049
$ws.Range("K21").Value2 = 1
$ws.Range("L21").Value2 = 0.3333333333333333
$ws.Range("M21").Value2 = 0.006136
$ws.Range("N21").Value2 = 0.018408
$ws.Range("O21").Value2 = 0.0005985839431548165
$ws.Range("P21").Value2 = 0.0006694856249710864
$ws.Range("Q21").Value2 = 0.004896260061333332
$ws.Range("R21").Value2 = 0.044066340552
$ws.Range("S21").Value2 = 0.00002455109961269137
$ws.Range("T21").Value2 = 0.00003555507702961083

# Row 22
$ws.Range("E22").Value2 = 2
$ws.Range("F22").Value2 = 0.6666666666666666
$ws.Range("G22").Value2 = 0.1455233333333333
$ws.Range("H22").Value2 = 0.43657
$ws.Range("I22").Value2 = 0.007479962034613753
$ws.Range("J22").Value2 = 0.009685317797424247
$ws.Range("M22").Value2 = 0.3250655
$ws.Range("N22").Value2 = 0.650131
$ws.Range("O22").Value2 = 0.03171104771407953
$ws.Range("P22").Value2 = 0.02364479350543662
$ws.Range("Q22").Value2 = 0.04730461511166667
$ws.Range("R22").Value2 = 0.28382769067
$ws.Range("S22").Value2 = 0.0002371974329791401
$ws.Range("T22").Value2 = 0.0002290073393546266

# Row 23
$ws.Range("E23").Value2 = 2
$ws.Range("F23").Value2 = 0.6666666666666666
$ws.Range("G23").Value2 = 0.1455233333333333
$ws.Range("H23").Value2 = 0.43657
$ws.Range("I23").Value2 = 0.007479962034613753
$ws.Range("J23").Value2 = 0.009685317797424247
$ws.Range("O23").Value2 = 0.6816872822276142
$ws.Range("P23").Value2 = 0.762432473166021
$ws.Range("Q23").Value2 = 1.016899687548889
$ws.Range("R23").Value2 = 9.152097187940001
$ws.Range("S23").Value2 = 0.005098994990541585
$ws.Range("T23").Value2 = 0.007384400801689048

# Row 24
$ws.Range("E24").Value2 = 2
$ws.Range("F24").Value2 = 0.6666666666666666
$ws.Range("G24").Value2 = 0.1455233333333333
$ws.Range("H24").Value2 = 0.43657
$ws.Range("I24").Value2 = 0.007479962034613753
$ws.Range("J24").Value2 = 0.009685317797424247
$ws.Range("M24").Value2 = 2.9317775
$ws.Range("N24").Value2 = 5.863555
$ws.Range("O24").Value2 = 0.2860030861151516
$ws.Range("P24").Value2 = 0.2132532477035712
$ws.Range("Q24").Value2 = 0.4266420343916666
$ws.Range("R24").Value2 = 2.55985220635
$ws.Range("S24").Value2 = 0.002139292225923702
$ws.Range("T24").Value2 = 0.00206542547534192

# Row 25
$ws.Range("D25").Value2 = "Resolving-Mac"
$ws.Range("E25").Value2 = 2
$ws.Range("F25").Value2 = 0.6666666666666666
$ws.Range("G25").Value2 = 0.1455233333333333
$ws.Range("H25").Value2 = 0.43657
$ws.Range("I25").Value2 = 0.007479962034613753
$ws.Range("J25").Value2 = 0.009685317797424247
$ws.Range("K25").Value2 = 1
$ws.Range("L25").Value2 = 0.3333333333333333
$ws.Range("M25").Value2 = 0.006136
$ws.Range("N25").Value2 = 0.018408
$ws.Range("O25").Value2 = 0.0005985839431548165
$ws.Range("P25").Value2 = 0.0006694856249710864
$ws.Range("Q25").Value2 = 0.0008929311733333334
$ws.Range("R25").Value2 = 0.00803638056
$ws.Range("S25").Value2 = 0.000004477385169327424
$ws.Range("T25").Value2 = 0.000006484181038652158

